# Generate Report for Handoff
# The "6ea5906b-87f6-4695-b6e9-88e0df925214.md" file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with refreshed
# timestamps, and a new handback-version-mismatch error message recorded
# in the per-locale "Error Detail" column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e796d752d3858138bb4613e2ef5dbf02382cc6e/e2e/6ea5906b-87f6-4695-b6e9-88e0df925214.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99fe65b5c98212006242552b92cc2669f4474714/e2e/6ea5906b-87f6-4695-b6e9-88e0df925214.md."

# --- Overview sheet: row 3 is the 6ea5906b-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-19 11:09:34"

# --- zh-cn sheet: row 3 is the 6ea5906b-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-19 11:09:22"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the 6ea5906b-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-10-19 11:09:34"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
